# Insert a new data row at row 267 on the active sheet (shifting rows
# 267-326 down to 268-327) and populate it with the new record's values.
# This mirrors the diff: old row 267 (D=44424, M=13583, P=226, ... rest)
# moves down to become row 268, while the brand-new row 267 keeps the
# same J/K/L/N/O/Q values as the old row267 but has its own Fecha (D),
# Precio promedio ponderado (M) and Precio $/Kg (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 267:326 down to 268:327, inserting a blank row at 267.
$ws.Rows(267).Insert()

# Populate the newly inserted row 267 with the new record.
$ws.Range("A267").Value2 = 10
$ws.Range("B267").Value = "Vega Modelo de Temuco"
$ws.Range("C267").Value = "La Araucanía"
$ws.Range("D267").Value2 = 44785
$ws.Range("E267").Value2 = 9
$ws.Range("F267").Value2 = 100112001
$ws.Range("G267").Value = "Berenjena"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value2 = 120
$ws.Range("K267").Value2 = 13000
$ws.Range("L267").Value2 = 14000
$ws.Range("M267").Value2 = 13333
$ws.Range("N267").Value = "`$/caja 60 unidades"
$ws.Range("O267").Value = "Región de Arica y Parinacota"
$ws.Range("P267").Value2 = 222
$ws.Range("Q267").Value2 = 60
$ws.Range("R267").Value = "Hortaliza"
